# Replicates Python's str.title() semantics: a letter is upper-cased when it
# is the first letter of a run of letters (i.e. preceded by a non-letter or
# start-of-string); every other letter in that run is lower-cased. Non-letter
# characters (spaces, hyphens, punctuation, digits, accented marks that .NET
# doesn't classify as letters, etc.) are left untouched.
function PyTitle([string]$s) {
    $lower = $s.ToLower()
    $re = [regex]"\p{L}+"
    $wordMatches = $re.Matches($lower)
    $result = $lower
    $n = $wordMatches.Count
    for ($i = $n - 1; $i -ge 0; $i--) {
        $mm = $wordMatches[$i]
        $word = $mm.Value
        $capitalized = $word.Substring(0,1).ToUpper() + $word.Substring(1)
        $result = $result.Substring(0, $mm.Index) + $capitalized + $result.Substring($mm.Index + $mm.Length)
    }
    return $result
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header columns to the new machine-friendly field names.
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# Title-case every state (col A) and municipality (col B) label in the data
# rows, e.g. "Pabellón de Arteaga" -> "Pabellón De Arteaga", "TOTAL" -> "Total".
for ($r = 2; $r -le 1611; $r++) {
    $a = $ws.Cells.Item($r, 1).Value2
    if ($a -ne $null) {
        $ws.Cells.Item($r, 1).Value = (PyTitle $a)
    }
    $b = $ws.Cells.Item($r, 2).Value2
    if ($b -ne $null) {
        $ws.Cells.Item($r, 2).Value = (PyTitle $b)
    }
}

# Drop the trailing footnote rows (sample size / source / author / date) that
# used to sit below the grand-total row, shrinking the sheet to A1:D1611.
$ws.Rows("1613:1617").Delete()
